# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Almeria" / "Lugo" rows (row 47 and row 48): names, and their
# "Casos activos" (column C) values, which were mismatched and are now
# corrected by reordering the rows.
$ws.Range("A47").Value = "Lugo"
$ws.Range("A48").Value = "Almeria"

$ws.Range("C47").Value = 5
$ws.Range("C48").Value = 72

# Update the "last updated" timestamp string (A1) from 05:46 to 06:16.
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 06:16"
